$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 3
$ws.Range("F6").Value = -6
$ws.Range("F8").Value = -2
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -5
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = 1
$ws.Range("F18").Value = -4
